$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change C38 from "PAGAMENTO" to "SALÁRIO"
$ws.Cells.Item(38, 3).Value = "SALÁRIO"

# Add new rows 204-207
$newRows = @(
    @(204, "ADIANTAMENTO DE SALÁRIO", "SALÁRIO"),
    @(205, "RESCISÃO TRABALHISTA", "RESCISÃO TRABALHISTA"),
    @(206, "OPERADOR DE MAQUINAS", "PRESTADOR DE SERVIÇO"),
    @(207, "CENOGRAFO", "PRESTADOR DE SERVIÇO")
)

foreach ($row in $newRows) {
    $rowNum = $row[0]
    $ws.Cells.Item($rowNum, 1).Value = $rowNum
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
}

# Update the defined name range to cover the new rows
$n = $wb.Names.Item("Codigos_Despesas")
$n.RefersTo = "='Codigos_Despesas'!`$A`$1:`$D`$207"
